# "Generate Report for Handback"
# Updates the localization-status report after a handback: the status text
# moves from "Ready for handoff" to "Handed back: in sync with en-US" on the
# Overview sheet and on each language sheet's Status column, and the
# per-file handback metadata (target file, handback file/link, handback
# timestamp) gets filled in on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both files ---
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# --- zh-cn / de-de sheets: Status column (C) for both rows ---
$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack
$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack

# --- Handback metadata: Latest Target File (I), Latest Handback File (J),
#     Latest Handback DateTime (K) ---
$file1Md  = "0ef44f93-0eec-4a23-9946-1ef5b797c5e5.md"
$file2Md  = "9dedd2e1-1ddf-4400-8172-d496dfa26871.md"

$zhcn.Range("I2").Value = $file1Md
$zhcn.Range("J2").Value = "0ef44f93-0eec-4a23-9946-1ef5b797c5e5.b988925ff3302ade0120de66030767b81b2c48f4.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-14 17:04:19"

$zhcn.Range("I3").Value = $file2Md
$zhcn.Range("J3").Value = "9dedd2e1-1ddf-4400-8172-d496dfa26871.2cc8cf516475ba284f713af3abb53a1a9d61a8a8.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-14 17:04:19"

$dede.Range("I2").Value = $file1Md
$dede.Range("J2").Value = "0ef44f93-0eec-4a23-9946-1ef5b797c5e5.b988925ff3302ade0120de66030767b81b2c48f4.de-de.xlf"
$dede.Range("K2").Value = "2016-08-14 17:04:29"

$dede.Range("I3").Value = $file2Md
$dede.Range("J3").Value = "9dedd2e1-1ddf-4400-8172-d496dfa26871.2cc8cf516475ba284f713af3abb53a1a9d61a8a8.de-de.xlf"
$dede.Range("K3").Value = "2016-08-14 17:04:29"

# --- Widen the columns that now hold the longer "Handed back: in sync with
#     en-US" status text and the handback file links ---
$overview.Columns.Item(5).ColumnWidth = 29.17   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.17   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth  = 29.17   # C: Status
$zhcn.Columns.Item(9).ColumnWidth  = 39.17   # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.17   # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth  = 29.17   # C: Status
$dede.Columns.Item(9).ColumnWidth  = 39.17   # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.17   # J: Latest Handback File

# --- Link the newly-populated "Latest Target File" cells to the source doc,
#     same as column A's existing handoff-file hyperlinks ---
$url1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/3626ec254d82a8ebceef698575800c108d20d53d/e2e/0ef44f93-0eec-4a23-9946-1ef5b797c5e5.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/3626ec254d82a8ebceef698575800c108d20d53d/e2e/9dedd2e1-1ddf-4400-8172-d496dfa26871.md"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $file1Md) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $file2Md) | Out-Null

$dede.Hyperlinks.Add($dede.Range("I2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $file1Md) | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $file2Md) | Out-Null
